# Update "想去人数" (interest count) figures in column F across sheets
# 展览 (Exhibitions), 演出 (Shows), and 全部类型 (All Types, an aggregate sheet).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 244
$wsExpo.Range("F3").Value  = 1434
$wsExpo.Range("F5").Value  = 896
$wsExpo.Range("F7").Value  = 1257
$wsExpo.Range("F8").Value  = 1605
$wsExpo.Range("F9").Value  = 160
$wsExpo.Range("F10").Value = 56
$wsExpo.Range("F11").Value = 2285
$wsExpo.Range("F12").Value = 461
$wsExpo.Range("F13").Value = 128
$wsExpo.Range("F17").Value = 85
$wsExpo.Range("F18").Value = 6284
$wsExpo.Range("F20").Value = 6192
$wsExpo.Range("F21").Value = 10186
$wsExpo.Range("F24").Value = 185
$wsExpo.Range("F26").Value = 505
$wsExpo.Range("F28").Value = 155
$wsExpo.Range("F30").Value = 105
$wsExpo.Range("F31").Value = 395

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 12
$wsShow.Range("F8").Value = 1161

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 244
$wsAll.Range("F5").Value  = 1434
$wsAll.Range("F8").Value  = 896
$wsAll.Range("F10").Value = 1257
$wsAll.Range("F11").Value = 12
$wsAll.Range("F12").Value = 1605
$wsAll.Range("F14").Value = 160
$wsAll.Range("F15").Value = 2285
$wsAll.Range("F17").Value = 461
$wsAll.Range("F18").Value = 128
$wsAll.Range("F23").Value = 85
$wsAll.Range("F24").Value = 6284
$wsAll.Range("F26").Value = 6192
$wsAll.Range("F27").Value = 10186
$wsAll.Range("F31").Value = 185
$wsAll.Range("F34").Value = 505
$wsAll.Range("F39").Value = 155
$wsAll.Range("F46").Value = 395
